# Generate Report for Handback
# Updates the timestamp values recorded on the handback-status report:
#  - Overview sheet's "Latest HO Xliff Generate Date" (G2) and the de-de
#    sheet's "Correspond Handoff Datetime" (H2) shared the same timestamp
#    string, so both move from 01:02:17 -> 01:02:58.
#  - zh-cn sheet's handoff/handback datetimes (H2/K2) advance to 01:02:54
#    and 01:03:13 respectively.
#  - de-de sheet's "Correspond Handback DateTime" (K2) advances to 01:03:20.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-08-28 01:02:58"

$wsZhCn.Range("H2").Value = "2016-08-28 01:02:54"
$wsZhCn.Range("K2").Value = "2016-08-28 01:03:13"

$wsDeDe.Range("H2").Value = "2016-08-28 01:02:58"
$wsDeDe.Range("K2").Value = "2016-08-28 01:03:20"
